$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("id>local_tomcat")
Write-Host "found at $($rng.Start) $($rng.End)"
$rng.Bold = 1
$rng.Bold = 0
